$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (B), shifting all subsequent
# columns (C:I -> B:H) one position to the left.
$ws.Columns("B").Delete()

# Append ".global" to each of the remaining header labels (B1:H1),
# leaving the "Country" header in A1 untouched.
$ws.Range("B1").Value = "Daily.global"
$ws.Range("C1").Value = "4 to 6 days per week.global"
$ws.Range("D1").Value = "2 to 3 days per week.global"
$ws.Range("E1").Value = "Once a week or less.global"
$ws.Range("F1").Value = "Not used in the last 30 days.global"
$ws.Range("G1").Value = "Not known / missing.global"
$ws.Range("H1").Value = "Total.global"
